$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.590.96'
$ws.Range("E2").Value = '  +2.97%  '

$ws.Range("D3").Value = '2.547.78'
$ws.Range("E3").Value = '  +1.52%  '

$ws.Range("D5").Value = '321.12'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").Value = '108.38'
$ws.Range("E6").Value = '  -0.40%  '

$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("D8").Value = '0.998'

$ws.Range("D9").Value = '0.553'
$ws.Range("E9").Value = '  +1.78%  '

$ws.Range("D10").Value = '40.02'
$ws.Range("E10").Value = '  +0.48%  '

$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").Value = '0.0813'
$ws.Range("E12").Value = '  -0.56%  '

$ws.Range("E13").Value = '  +0.88%  '

$ws.Range("D14").Value = '7.22'
$ws.Range("E14").Value = '  +0.31%  '

$ws.Range("D15").Value = '2.944.83'
$ws.Range("E15").Value = '  +1.49%  '

$ws.Range("D16").Value = '2.564.41'
$ws.Range("E16").Value = '  +2.43%  '

$ws.Range("D17").Value = '0.855'
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").Value = '49.347.21'
$ws.Range("E18").Value = '  +2.79%  '

$ws.Range("D19").Value = '13.18'
$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("D20").Value = '2.94'
$ws.Range("E20").Value = '  +7.71%  '

$ws.Range("E21").Value = '  +1.26%  '

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").Value = '283.39'
$ws.Range("E23").Value = '  +3.24%  '

$ws.Range("D24").Value = '71.57'
$ws.Range("E24").Value = '  -0.95%  '

$ws.Range("E25").Value = '  -1.58%  '

$ws.Range("D26").Value = '26.27'
$ws.Range("E26").Value = '  +1.46%  '

$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '0.145'
$ws.Range("E28").Value = '  +3.39%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  -7.62%  '

$ws.Range("D30").Value = '9.79'
$ws.Range("E30").Value = '  -2.62%  '

$ws.Range("D31").Value = '35.21'
$ws.Range("E31").Value = '  -0.73%  '

$ws.Range("D32").Value = '49.58'
$ws.Range("E32").Value = '  +0.17%  '

$ws.Range("D33").Value = '19.64'
$ws.Range("E33").Value = '  +1.58%  '

$ws.Range("D34").Value = '5.36'
$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("E35").Value = '  -0.15%  '

$ws.Range("D36").Value = '0.0782'
$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("E37").Value = '  +2.38%  '

$ws.Range("E38").Value = '  +0.45%  '

$ws.Range("E39").Value = '  -0.17%  '

$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("E41").Value = '  +0.48%  '

$ws.Range("D42").Value = '22.14'
$ws.Range("E42").Value = '  +2.29%  '

$ws.Range("D43").Value = '120.08'
$ws.Range("E43").Value = '  -2.00%  '

$ws.Range("D44").Value = '0.0309'
$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("D45").Value = '3.27'
$ws.Range("E45").Value = '  +4.69%  '

$ws.Range("D46").Value = '2.012.08'
$ws.Range("E46").Value = '  -0.35%  '

$ws.Range("D47").Value = '1.98'
$ws.Range("E47").Value = '  +6.79%  '

$ws.Range("D48").Value = '2.12'
$ws.Range("E48").Value = '  +6.66%  '

$ws.Range("D49").Value = '9.02'
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("D50").Value = '5.30'
$ws.Range("E50").Value = '  +2.20%  '

$ws.Range("D51").Value = '81.13'
$ws.Range("E51").Value = '  +2.07%  '
